# Add "Pay Later" column to the order sample sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in O1 and value in O2
$ws.Range("O1").Value = "Pay Later"
$ws.Range("O2").Value = "Yes"

# Move the active selection to O3, matching the saved selection state
$ws.Range("O3").Select()
